# Commit: add color blindness toggle
#
# Simplifies the frequency / fat-type answer labels used in the dataset so
# a colour-blind-friendly legend (single short words) can be used instead
# of the long verbose labels. This collapses a couple of redundant
# "matiere_grasse" choices ("Margarine" / "Margarine et huile vegetale")
# into the existing "Beurre" choice.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$whole = [Microsoft.Office.Interop.Excel.XlLookAt]::xlWhole

# Shorten the consumption-frequency labels (frequence_consommation_* columns)
$ws.Cells.Replace("rarement 0-1 fois/semaine", "rarement", $whole)
$ws.Cells.Replace("occasionnellement 2-3 fois/semaine", "occasionnellement", $whole)
$ws.Cells.Replace("Fréquente >4 fois/semaine", "Fréquente", $whole)

# Collapse the butter/margarine fat-type answers into a single "Beurre" choice
$ws.Cells.Replace("Beurre ou beurre allégé", "Beurre", $whole)
$ws.Cells.Replace("Margarine et huile végétale", "Beurre", $whole)
$ws.Cells.Replace("Margarine", "Beurre", $whole)
